$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.732.56'
$ws.Range("E2").Value = '  +2.07%  '

$ws.Range("D3").Value = '2.115.85'
$ws.Range("E3").Value = '  +10.36%  '

$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '333.80'
$ws.Range("E5").Value = '  +4.28%  '

$ws.Range("E6").Value = '  -0.08%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5233'
$ws.Range("E7").Value = '  +3.36%  '

$ws.Range("E8").Value = '  +8.55%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.09063'
$ws.Range("E9").Value = '  +8.52%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '46.44'
$ws.Range("E10").Value = '  +9.54%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.182'
$ws.Range("E11").Value = '  +6.52%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '25.26'
$ws.Range("E12").Value = '  +5.13%  '

$ws.Range("D13").Value = '2.111.99'
$ws.Range("E13").Value = '  +10.02%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.818'
$ws.Range("E14").Value = '  +6.31%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '98.39'
$ws.Range("E16").Value = '  +6.31%  '

$ws.Range("E17").Value = '  -0.25%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001140'
$ws.Range("E18").Value = '  +3.93%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06655'
$ws.Range("E19").Value = '  +2.18%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.29'
$ws.Range("E20").Value = '  +4.00%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.410'
$ws.Range("E21").Value = '  +7.83%  '

$ws.Range("E22").Value = '  -0.05%  '

$ws.Range("D23").Value = '30.856.13'
$ws.Range("E23").Value = '  +2.42%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.08'
$ws.Range("E24").Value = '  +6.29%  '

$ws.Range("D25").Value = '2.361.10'
$ws.Range("E25").Value = '  +10.52%  '

$ws.Range("E26").Value = '  +3.06%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.97'
$ws.Range("E27").Value = '  +5.23%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.554'
$ws.Range("E28").Value = '  +12.38%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '163.41'
$ws.Range("E29").Value = '  +0.34%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '133.94'
$ws.Range("E30").Value = '  +3.92%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.187'
$ws.Range("E31").Value = '  +3.82%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.249'
$ws.Range("E33").Value = '  +5.04%  '

$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.915'
$ws.Range("E34").Value = '  +3.37%  '

$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.532'
$ws.Range("E35").Value = '  +28.00%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02591'
$ws.Range("E36").Value = '  +5.51%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.601'
$ws.Range("E37").Value = '  +4.41%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06784'
$ws.Range("E38").Value = '  +5.20%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '9.579'
$ws.Range("E39").Value = '  +11.09%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '12.76'
$ws.Range("E40").Value = '  +11.87%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2270'
$ws.Range("E41").Value = '  +5.26%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6807'
$ws.Range("E42").Value = '  +4.37%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.254'
$ws.Range("E43").Value = '  +3.64%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.24'
$ws.Range("E44").Value = '  +6.58%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9995'
$ws.Range("E45").Value = '  +0.04%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6356'
$ws.Range("E46").Value = '  +4.58%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.252'
$ws.Range("E47").Value = '  +2.85%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.293'
$ws.Range("E48").Value = '  +6.92%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.673'
$ws.Range("E49").Value = '  +1.38%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '83.13'
$ws.Range("E50").Value = '  +5.18%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '120.46'
$ws.Range("E51").Value = '  -1.54%  '

